# Auto-generated script to update FFXIV leve profit market data
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M8").Value = -10358
$ws.Range("H8").Value = 3499
$ws.Range("I8").Value = 3499
$ws.Range("K8").Value = 10497
$ws.Range("K33").Value = 164.28572
$ws.Range("I33").Value = 164.28572
$ws.Range("M33").Value = 64.71428
$ws.Range("H33").Value = 211
$ws.Range("K43").Value = 1049.4
$ws.Range("L43").Value = 2181.8572
$ws.Range("M43").Value = -980.4000000000001
$ws.Range("N43").Value = -2319.8572
$ws.Range("J43").Value = 2181.8572
$ws.Range("I43").Value = 1049.4
$ws.Range("H43").Value = 1710
$ws.Range("K51").Value = 10133.333
$ws.Range("M51").Value = -9649.333000000001
$ws.Range("H51").Value = 10066.167
$ws.Range("I51").Value = 10133.333
$ws.Range("J53").Value = 2500
$ws.Range("H53").Value = 4170
$ws.Range("L53").Value = 2500
$ws.Range("N53").Value = -3774
$ws.Range("K62").Value = 31254272
$ws.Range("I62").Value = 31254272
$ws.Range("H62").Value = 19235398
$ws.Range("M62").Value = -31253648
$ws.Range("M65").Value = -156268240
$ws.Range("I65").Value = 31254272
$ws.Range("K65").Value = 156271360
$ws.Range("H65").Value = 19235398
$ws.Range("H98").Value = 5918.4814
$ws.Range("K98").Value = 681.75
$ws.Range("J98").Value = 13535.546
$ws.Range("L98").Value = 13535.546
$ws.Range("M98").Value = 816.25
$ws.Range("I98").Value = 681.75
$ws.Range("N98").Value = -16531.546
$ws.Range("H99").Value = 1133.3334
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("K111").Value = 7842
$ws.Range("H111").Value = 2614
$ws.Range("M111").Value = -4775
$ws.Range("I111").Value = 2614
$ws.Range("N112").Value = -28403.999
$ws.Range("J112").Value = 8729.333000000001
$ws.Range("H112").Value = 8729.333000000001
$ws.Range("L112").Value = 26187.999
$ws.Range("H113").Value = 6330.893
$ws.Range("M113").Value = -2719.4443
$ws.Range("I113").Value = 5973.4443
$ws.Range("K113").Value = 5973.4443
$ws.Range("L122").Value = 40606.638
$ws.Range("N122").Value = -45506.638
$ws.Range("J122").Value = 13535.546
$ws.Range("M122").Value = 404.75
$ws.Range("I122").Value = 681.75
$ws.Range("K122").Value = 2045.25
$ws.Range("H122").Value = 5918.4814
$ws.Range("K132").Value = 8550.9645
$ws.Range("N132").Value = -20022.5
$ws.Range("H132").Value = 2924.0173
$ws.Range("M132").Value = -6020.9645
$ws.Range("J132").Value = 4987.5
$ws.Range("I132").Value = 2850.3215
$ws.Range("L132").Value = 14962.5
$ws.Range("H135").Value = 1613.5518
$ws.Range("M135").Value = -9650.249400000001
$ws.Range("N135").Value = -30808.2
$ws.Range("I135").Value = 1353.9166
$ws.Range("J135").Value = 2859.8
$ws.Range("K135").Value = 12185.2494
$ws.Range("L135").Value = 25738.2
$ws.Range("N99").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 2034.2742
$ws.Range("I32").Value = 2034.2742
$ws.Range("H32").Value = 2595.5
$ws.Range("M32").Value = -1747.2742
$ws.Range("L122").Value = 10497.4284
$ws.Range("N122").Value = -15397.4284
$ws.Range("J122").Value = 3499.1428
$ws.Range("M122").Value = -622
$ws.Range("I122").Value = 1024
$ws.Range("K122").Value = 3072
$ws.Range("H122").Value = 2106.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J20").Value = 931.4
$ws.Range("H20").Value = 913.2632
$ws.Range("N20").Value = -1425.4
$ws.Range("L20").Value = 931.4
$ws.Range("N86").Value = -13187.23
$ws.Range("M86").Value = -2565.1
$ws.Range("K86").Value = 3688.1
$ws.Range("J86").Value = 10941.23
$ws.Range("L86").Value = 10941.23
$ws.Range("I86").Value = 3688.1
$ws.Range("H86").Value = 6545.394
$ws.Range("I89").Value = 3688.1
$ws.Range("K89").Value = 18440.5
$ws.Range("J89").Value = 10941.23
$ws.Range("N89").Value = -65938.14999999999
$ws.Range("L89").Value = 54706.14999999999
$ws.Range("H89").Value = 6545.394
$ws.Range("M89").Value = -12824.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I99").Value = 2028.5238
$ws.Range("H99").Value = 2081.7727
$ws.Range("K99").Value = 2028.5238
$ws.Range("M99").Value = -530.5237999999999
$ws.Range("L105").Value = 2500
$ws.Range("K105").Value = 896.75
$ws.Range("N105").Value = -5994
$ws.Range("I105").Value = 896.75
$ws.Range("H105").Value = 1217.4
$ws.Range("J105").Value = 2500
$ws.Range("M105").Value = 850.25
$ws.Range("H126").Value = 2081.7727
$ws.Range("M126").Value = -3615.5714
$ws.Range("I126").Value = 2028.5238
$ws.Range("K126").Value = 6085.5714
$ws.Range("K132").Value = 10515.4614
$ws.Range("N132").Value = -17395.1432
$ws.Range("H132").Value = 3633.818
$ws.Range("M132").Value = -7985.4614
$ws.Range("J132").Value = 4111.7144
$ws.Range("I132").Value = 3505.1538
$ws.Range("L132").Value = 12335.1432
$ws.Range("M134").Value = -1318.2
$ws.Range("I134").Value = 1284.4
$ws.Range("K134").Value = 3853.2
$ws.Range("H134").Value = 1848.9546
$ws.Range("H141").Value = 234059.1
$ws.Range("L141").Value = 273961.12
$ws.Range("N141").Value = -284321.12
$ws.Range("J141").Value = 273961.12

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K5").Value = 1922.0001
$ws.Range("N5").Value = -2736.5
$ws.Range("H5").Value = 663.82355
$ws.Range("J5").Value = 837.5
$ws.Range("M5").Value = -1810.0001
$ws.Range("I5").Value = 640.6667
$ws.Range("L5").Value = 2512.5
$ws.Range("M8").Value = -2067.66662
$ws.Range("H8").Value = 735.55554
$ws.Range("I8").Value = 735.55554
$ws.Range("K8").Value = 2206.66662
$ws.Range("L62").Value = 32995.5
$ws.Range("H62").Value = 10998.5
$ws.Range("J62").Value = 10998.5
$ws.Range("N62").Value = -34367.5
$ws.Range("J65").Value = 10998.5
$ws.Range("L65").Value = 98986.5
$ws.Range("N65").Value = -105850.5
$ws.Range("H65").Value = 10998.5
$ws.Range("I107").Value = 490.5
$ws.Range("N107").Value = -5726.1429
$ws.Range("M107").Value = 448.5
$ws.Range("H107").Value = 619.5
$ws.Range("J107").Value = 628.7143
$ws.Range("K107").Value = 1471.5
$ws.Range("L107").Value = 1886.1429
$ws.Range("H135").Value = 663.82355
$ws.Range("M135").Value = -3231.0003
$ws.Range("N135").Value = -12607.5
$ws.Range("I135").Value = 640.6667
$ws.Range("J135").Value = 837.5
$ws.Range("K135").Value = 5766.0003
$ws.Range("L135").Value = 7537.5
$ws.Range("H138").Value = 38473196
$ws.Range("J138").Value = 16746.6
$ws.Range("N138").Value = -60519.8
$ws.Range("L138").Value = 50239.8
$ws.Range("K140").Value = 6926.25
$ws.Range("M140").Value = -1746.25
$ws.Range("N140").Value = -58357.624
$ws.Range("L140").Value = 47997.624
$ws.Range("H140").Value = 10523.025
$ws.Range("I140").Value = 2308.75
$ws.Range("J140").Value = 15999.208

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 39000
$ws.Range("N47").Value = -40136
$ws.Range("L47").Value = 39000
$ws.Range("J47").Value = 39000

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M40").Value = -2420.111
$ws.Range("N40").Value = -3471
$ws.Range("L40").Value = 3199
$ws.Range("J40").Value = 3199
$ws.Range("I40").Value = 2556.111
$ws.Range("K40").Value = 2556.111
$ws.Range("H40").Value = 2673
$ws.Range("H61").Value = 1573.1666
$ws.Range("K61").Value = 1487.8
$ws.Range("I61").Value = 1487.8
$ws.Range("M61").Value = -1285.8
$ws.Range("H113").Value = 1573.1666
$ws.Range("M113").Value = 682.2
$ws.Range("I113").Value = 1487.8
$ws.Range("K113").Value = 1487.8
$ws.Range("H118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L122").Value = 11100
$ws.Range("N122").Value = -16000
$ws.Range("J122").Value = 3700
$ws.Range("M122").Value = -8496.000100000001
$ws.Range("I122").Value = 3648.6667
$ws.Range("K122").Value = 10946.0001
$ws.Range("H122").Value = 3674.3333
$ws.Range("N134").Value = -263765
$ws.Range("L134").Value = 253625
$ws.Range("H134").Value = 253625
$ws.Range("J134").Value = 253625
$ws.Range("I136").Value = 2330.5557
$ws.Range("N136").Value = -14077.5
$ws.Range("K136").Value = 6991.6671
$ws.Range("J136").Value = 2992.5
$ws.Range("H136").Value = 2474.4565
$ws.Range("L136").Value = 8977.5
$ws.Range("M136").Value = -4441.6671
$ws.Range("N118").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N100").Value = -5803
$ws.Range("H100").Value = 1435.1428
$ws.Range("I100").Value = 1065
$ws.Range("L100").Value = 4721
$ws.Range("J100").Value = 2360.5
$ws.Range("M100").Value = -1589
$ws.Range("K100").Value = 2130
$ws.Range("I107").Value = 410.375
$ws.Range("N107").Value = -6290.4
$ws.Range("M107").Value = 688.875
$ws.Range("H107").Value = 636.1667
$ws.Range("J107").Value = 816.8
$ws.Range("K107").Value = 1231.125
$ws.Range("L107").Value = 2450.4
$ws.Range("N112").Value = -26416.334
$ws.Range("J112").Value = 23462.334
$ws.Range("H112").Value = 23462.334
$ws.Range("L112").Value = 23462.334
$ws.Range("L113").Value = 3256.8
$ws.Range("N113").Value = -7596.799999999999
$ws.Range("H113").Value = 613.4074000000001
$ws.Range("M113").Value = 651.7273
$ws.Range("I113").Value = 506.0909
$ws.Range("K113").Value = 1518.2727
$ws.Range("J113").Value = 1085.6
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("N116").Value = -39178
$ws.Range("L116").Value = 30000
$ws.Range("L122").Value = 154422
$ws.Range("N122").Value = -159322
$ws.Range("J122").Value = 51474
$ws.Range("M122").Value = -13918.9171
$ws.Range("I122").Value = 5456.3057
$ws.Range("K122").Value = 16368.9171
$ws.Range("H122").Value = 7878.2896
$ws.Range("I136").Value = 2385.0378
$ws.Range("N136").Value = -57245.25
$ws.Range("K136").Value = 7155.1134
$ws.Range("J136").Value = 17381.75
$ws.Range("H136").Value = 3437.4385
$ws.Range("L136").Value = 52145.25
$ws.Range("M136").Value = -4605.1134
